$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "77÷6="
$t.Cell(1,2).Range.Text = "35÷8="
$t.Cell(1,3).Range.Text = "26÷5="
$t.Cell(1,4).Range.Text = "83÷6="
$t.Cell(1,5).Range.Text = "40÷8="
$t.Cell(5,1).Range.Text = "39÷3="
$t.Cell(5,2).Range.Text = "85÷5="
$t.Cell(5,3).Range.Text = "76÷7="
$t.Cell(5,4).Range.Text = "13÷3="
$t.Cell(5,5).Range.Text = "17÷4="
$t.Cell(9,1).Range.Text = "15÷9="
$t.Cell(9,2).Range.Text = "51÷8="
$t.Cell(9,3).Range.Text = "94÷9="
$t.Cell(9,4).Range.Text = "61÷9="
$t.Cell(9,5).Range.Text = "70÷7="
$t.Cell(13,1).Range.Text = "18÷2="
$t.Cell(13,2).Range.Text = "90÷7="
$t.Cell(13,3).Range.Text = "75÷3="
$t.Cell(13,4).Range.Text = "92÷7="
$t.Cell(13,5).Range.Text = "71÷2="
$t.Cell(17,1).Range.Text = "82÷7="
$t.Cell(17,2).Range.Text = "55÷9="
$t.Cell(17,3).Range.Text = "28÷9="
$t.Cell(17,4).Range.Text = "62÷5="
$t.Cell(17,5).Range.Text = "62÷8="
